$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4549
$ws1.Range("F3").Value = 2499
$ws1.Range("F5").Value = 31
$ws1.Range("F12").Value = 1702
$ws1.Range("F14").Value = 3738
$ws1.Range("F15").Value = 24
$ws1.Range("F16").Value = 249

# Sheet "全部类型" (all types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4549
$ws4.Range("F3").Value = 2499
$ws4.Range("F5").Value = 31
$ws4.Range("F16").Value = 1702
$ws4.Range("F18").Value = 3738
$ws4.Range("F19").Value = 24
$ws4.Range("F20").Value = 249
